# Apply the edits described by the commit: add the batter-pitcher
# interaction row (row 35) on Sheet2, and update the active
# sheet/selection state so Sheet1 is the active tab while Sheet2's
# selection moves to B29 (no longer scrolled/tabSelected).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- New "batter-pitcher interaction" data in row 35 of Sheet2 ---
$ws2.Range("D35").Formula = "=993/4083"
$ws2.Range("E35").Formula = "=346/4083"
$ws2.Range("F35").Value = 0.83
$ws2.Range("G35").Value = 0.122
$ws2.Range("H35").Value = 0.272
$ws2.Range("I35").Formula = "=254/975"

# --- Update Sheet2's view: selection moves to B29, scroll resets ---
$ws2.Range("B29").Select() | Out-Null

# --- Sheet1 becomes the active/selected tab (was Sheet2 before) ---
$ws1.Activate() | Out-Null
$ws1.Range("D2").Select() | Out-Null
